# Edit script for "Examen Mora.docx" — applies the changes described by the
# commit diff using Word COM-interop (Find/Replace + paragraph-level edits).

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $ok = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "WARNING: replace failed for:" $old
    }
}

# --- Paragraph 1 ("The topic that I am choosing...") ---------------------
Replace-Text "The topic that I am choosing is Text Summarization." "The topic that I chose is Text Summarization."

# --- Paragraph 3 ("For this topic, ... For project, ...") ----------------
Replace-Text "that it would be preferably to exist" "that would be preferable to exist"
Replace-Text "text/paragraph/article etc. For project," "text/paragraph/article etc. For the project,"

# --- Paragraph 4 (long "contribute ..." paragraph) ------------------------
Replace-Text "contribute in the transformation process because them being very common and have little to no meaning. Once this pre-processing step is done, the next one counting all of the words in the paragraph and their appearance. This step is done because we can see the most and least words that appearance and outline the topics." "contribute to the transformation process because of them being very common and having little to no meaning. Once this pre-processing step is done, the next one is counting all of the words in the paragraph and their appearance. This step is done because we can see the most and least words that appear and outline the topics."

Replace-Text "we calculate the weight occurrence frequency of all the words. To achieve this, we calculate by the most recurrent word in the paragraph. For example, if the most appeared word is “forest” and the second one is “tree”, we have:" "we calculate the weight occurrence frequency of all the words."

# Append a trailing red space run after "...the words." in paragraph 4.
$p4 = $d.Paragraphs(4)
$insertPoint = $d.Range($p4.Range.End - 1, $p4.Range.End - 1)
$insertPoint.InsertAfter(" ")
$insertPoint.Font.Color = 255
$insertPoint.LanguageID = "en-US"

# --- Paragraph 5 (previously empty ListParagraph) gets new content -------
$p5 = $d.Paragraphs(5)
$p5.Range.Text = "Practically, very similar to the project, the following step is eliminating the most and least common topics using the 5-95% rule which was the result of multiple simulations. By eliminating topics that have the most amount and the least amount of appearance, not only we reduced the number of articles and topics, but also the number of words. This step is exceptionally useful because of its optimized results. "
$p5.Range.LanguageID = "en-US"

# --- Paragraphs 6,7,8 (the "Word/Freq/Weight", "Forest", "Tree" rows) -----
# are removed entirely (the little results table made of tabbed paragraphs).
$d.Paragraphs(8).Range.Delete()   # Tree   2   0.67
$d.Paragraphs(7).Range.Delete()   # Forest 3   1
$d.Paragraphs(6).Range.Delete()   # Word Freq Weight

# --- Paragraph 6 (was "The next step required...") becomes a ListParagraph
#     bullet item with all-new text.
$p6 = $d.Paragraphs(6)
$p6.Style = "List Paragraph"
$p6.Range.Text = "Next, we determine the range of values and where each word is located by calculating its entropy that is utilized to calculate the information gain necessary in the proceeding step. The latter is useful in depicting the connection between words and topics. "
$p6.Range.LanguageID = "en-US"

# --- Paragraph 7 (previously empty ListParagraph+numPr) gets new content -
$p7 = $d.Paragraphs(7)
$p7.Range.Text = "The 2 approaches to summarize a text are extraction-based and abstraction-based summarization. For the former, the extraction is represented by pulling only the most information out of a piece of text and combining it to create a summary. For the latter, is used to paraphrase and shorten the original text. The major difference between these is observed in the complexity and text quality making abstraction-based summarization technique more complicated and sophisticated."
$p7.Range.LanguageID = "en-US"

# --- Insert a brand-new bullet paragraph after it --------------------------
$p7.Range.InsertParagraphAfter()
$p8 = $d.Paragraphs(8)
$p8.Range.Text = "A learning algorithm that could improve the quality of the text is inverse reinforcement learning IRL that focuses on estimating the reward function of an agent given a set of observations of that agent" + [char]0x2019 + "s behavior. Generally, IRL provides advantages in situations where the reward function is not explicitly known or where it is difficult to define or interact with the environment directly. These situations are exactly what we observe in summarization."
$p8start = $p8.Range.Start
$boundary = $p8start + 213
$d.Range($p8start, $boundary).LanguageID = "en-US"

Write-Host "done"
